$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell D10: "T50I  T50I 3UTR " -> "  T50I 3UTR "
$ws.Range("D10").Value = "  T50I 3UTR "

# Update cells that held "Q61R Q61R   Missense_Mutation" -> " Q61R   Missense_Mutation"
$ws.Range("D11").Value = " Q61R   Missense_Mutation"
$ws.Range("D12").Value = " Q61R   Missense_Mutation"
$ws.Range("D19").Value = " Q61R   Missense_Mutation"
$ws.Range("D28").Value = " Q61R   Missense_Mutation"
$ws.Range("D33").Value = " Q61R   Missense_Mutation"
$ws.Range("D44").Value = " Q61R   Missense_Mutation"
